# The template used a "manual" Word field (fldChar begin / instrText* /
# fldChar end) to hold the M2Doc script:
#
#   { m:'anydsl class diagram'.representationByName().asImage('PNG').fit(400, 400) }
#
# The parser was switched to TokenIteratorFieldRewriterSplit, which expects
# the M2Doc tag to be stored as plain literal text (a run per original
# instrText fragment) delimited by the usual "{" / "}" tag markers, instead
# of being wrapped in a real Word field. This rewrites that single
# paragraph: the fldChar begin/end runs are dropped, every w:instrText run
# becomes a w:t run with identical text, the first run gains a leading "{"
# and the last run gains a trailing "}".

$d = $word.ActiveDocument

# Locate the (single) legacy field and the paragraph range that contains it.
$field = $d.Fields.Item(1)
$fieldStart = $field.Code.Start

$target = $null
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($fieldStart -ge $r.Start -and $fieldStart -lt $r.End) {
        $target = $r
    }
}
if ($target -eq $null) {
    throw "Could not locate the paragraph containing the legacy field."
}

# Rebuild the paragraph as plain-text runs (same rPr / rsid attributes as
# the original instrText runs), wrapping the whole M2Doc expression in the
# "{" / "}" tag delimiters that used to be implied by the field braces.
$newParagraphXml = '<w:p w14:paraId="70197535" w14:textId="4ACFCB7D" w:rsidR="00A10D75" w:rsidRDefault="00474E78" w:rsidP="00F65375">' + `
  '<w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{m:</w:t></w:r>' + `
  '<w:r w:rsidR="00C53443"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>' + "'" + '</w:t></w:r>' + `
  '<w:r w:rsidR="008E7B0B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>anydsl class diagram</w:t></w:r>' + `
  '<w:r w:rsidR="00C53443"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>' + "'." + '</w:t></w:r>' + `
  '<w:r w:rsidR="007B65B7"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>r</w:t></w:r>' + `
  '<w:r w:rsidR="009B6BD9" w:rsidRPr="009B6BD9"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>epresentation</w:t></w:r>' + `
  '<w:r w:rsidR="007B65B7"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>By</w:t></w:r>' + `
  '<w:r w:rsidR="009B6BD9" w:rsidRPr="009B6BD9"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Name</w:t></w:r>' + `
  '<w:r w:rsidR="00C53443"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>()</w:t></w:r>' + `
  '<w:r w:rsidR="002D294B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r>' + `
  '<w:r w:rsidR="00D75A30"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>asImage(</w:t></w:r>' + `
  '<w:r w:rsidR="00A04964"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>' + "'" + '</w:t></w:r>' + `
  '<w:r w:rsidR="009C3812"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>PNG</w:t></w:r>' + `
  '<w:r w:rsidR="00A04964"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>' + "'" + '</w:t></w:r>' + `
  '<w:r w:rsidR="00D75A30"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>)</w:t></w:r>' + `
  '<w:r w:rsidR="000D38C7"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.fit(400, 400)}</w:t></w:r>' + `
  '</w:p>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
  '<w:body>' + $newParagraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($packageXml)

Write-Output "Field unwrapped into plain-text tag runs."
